$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.228.48"
$ws.Range("E2").Value = "  -1.52%  "

$ws.Range("D3").Value = "3.954.36"
$ws.Range("E3").Value = "  -1.67%  "

$ws.Range("E4").Value = "  -0.10%  "

$ws.Range("D5").Value = "'543.13"
$ws.Range("E5").Value = "  +5.52%  "

$ws.Range("D6").Value = "'148.90"
$ws.Range("E6").Value = "  +1.65%  "

$ws.Range("D7").Value = "'0.695"
$ws.Range("E7").Value = "  -2.72%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").Value = "'0.740"
$ws.Range("E9").Value = "  -3.21%  "

$ws.Range("D10").Value = "'0.168"
$ws.Range("E10").Value = "  -2.78%  "

$ws.Range("D11").Value = "'53.55"
$ws.Range("E11").Value = "  +15.18%  "

$ws.Range("D12").Value = "'0.0000320"
$ws.Range("E12").Value = "  -0.66%  "

$ws.Range("D13").Value = "'10.54"

$ws.Range("D14").Value = "4.609.20"
$ws.Range("E14").Value = "  -1.22%  "

$ws.Range("D15").Value = "3.960.33"
$ws.Range("E15").Value = "  -1.63%  "

$ws.Range("D16").Value = "'13.94"
$ws.Range("E16").Value = "  -0.60%  "

$ws.Range("D17").Value = "'20.34"
$ws.Range("E17").Value = "  -3.30%  "

$ws.Range("E18").Value = "  -0.53%  "

$ws.Range("E19").Value = "  -2.13%  "

$ws.Range("D20").Value = "71.253.03"
$ws.Range("E20").Value = "  -1.25%  "

$ws.Range("D21").Value = "'426.48"
$ws.Range("E21").Value = "  -1.95%  "

$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").Value = "'96.45"
$ws.Range("E22").Value = "  -5.74%  "

$ws.Range("B23").Value = "ImmutableX"
$ws.Range("C23").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D23").Value = "'3.53"
$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").Value = "'4.26"
$ws.Range("E24").Value = "  +7.72%  "

$ws.Range("D25").Value = "'14.21"
$ws.Range("E25").Value = "  -2.61%  "

$ws.Range("D26").Value = "'11.19"
$ws.Range("E26").Value = "  -3.12%  "

$ws.Range("D27").Value = "'10.51"
$ws.Range("E27").Value = "  -4.45%  "

$ws.Range("E28").Value = "  +1.18%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'3.64"
$ws.Range("E29").Value = "  +17.68%  "

$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "'36.37"
$ws.Range("E30").Value = "  -1.95%  "

$ws.Range("D31").Value = "'7.41"
$ws.Range("E31").Value = "  +9.05%  "

$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "'13.38"
$ws.Range("E32").Value = "  -1.16%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.130"
$ws.Range("E33").Value = "  +1.90%  "

$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").Value = "'678.19"
$ws.Range("E34").Value = "  +0.05%  "

$ws.Range("B35").Value = "InjectiveProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D35").Value = "'48.12"
$ws.Range("E35").Value = "  +17.89%  "

$ws.Range("D36").Value = "'65.48"
$ws.Range("E36").Value = "  -3.52%  "

$ws.Range("D37").Value = "'0.431"
$ws.Range("E37").Value = "  +0.18%  "

$ws.Range("D38").Value = "0.0₃0819"
$ws.Range("E38").Value = "  -5.33%  "

$ws.Range("D39").Value = "'0.149"
$ws.Range("E39").Value = "  -1.61%  "

$ws.Range("D40").Value = "'3.36"
$ws.Range("E40").Value = "  -4.00%  "

$ws.Range("D41").Value = "'3.35"
$ws.Range("E41").Value = "  +4.31%  "

$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.13%  "

$ws.Range("E43").Value = "  +0.26%  "

$ws.Range("D44").Value = "'0.0481"
$ws.Range("E44").Value = "  -0.86%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.148"
$ws.Range("E45").Value = "  -5.62%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "'2.69"
$ws.Range("E46").Value = "  -1.22%  "

$ws.Range("D47").Value = "'9.74"
$ws.Range("E47").Value = "  +8.17%  "

$ws.Range("D48").Value = "'3.37"
$ws.Range("E48").Value = "  -2.87%  "

$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'2.96"
$ws.Range("E49").Value = "  -3.05%  "

$ws.Range("B50").Value = "FLOKI"
$ws.Range("C50").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D50").Value = "'0.000273"
$ws.Range("E50").Value = "  +1.53%  "

$ws.Range("D51").Value = "'145.42"
$ws.Range("E51").Value = "  +2.15%  "
